# Fill in the previously-blank "Demand" value (column D) and the route / node
# sequence (columns F onward, up to the trailing -1 sentinel) for each
# solution row (7-21) on the active sheet. Column E already holds the -1
# sentinel and is left untouched; it is only restated here for clarity/safety.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    7  = @(125.0, -1.0, 12.0, 18.0, 41.0, 25.0, -1.0)
    8  = @(141.0, -1.0, 12.0, 33.0, 49.0, 30.0, 34.0, 9.0, 38.0, -1.0)
    9  = @(150.0, -1.0, 27.0, 48.0, 8.0, 31.0, 28.0, 2.0, 34.0, 39.0, -1.0)
    10 = @(144.0, -1.0, 12.0, 45.0, 33.0, 10.0, 49.0, 9.0, 38.0, 35.0, 3.0, -1.0)
    11 = @(81.0, -1.0, 46.0, 5.0, 39.0, 44.0, 42.0, 19.0, 17.0, -1.0)
    12 = @(92.0, -1.0, 6.0, 23.0, 7.0, 43.0, 24.0, 14.0, -1.0)
    13 = @(151.0, -1.0, 47.0, 18.0, 41.0, 40.0, 13.0, 25.0, -1.0)
    14 = @(153.0, -1.0, 12.0, 37.0, 41.0, 42.0, 44.0, 15.0, 5.0, 20.0, -1.0)
    15 = @(153.0, -1.0, 47.0, 4.0, 18.0, 41.0, 13.0, 25.0, -1.0)
    16 = @(146.0, -1.0, 27.0, 48.0, 26.0, 8.0, 20.0, 2.0, 34.0, -1.0)
    17 = @(124.0, -1.0, 1.0, 31.0, 28.0, 22.0, 20.0, 2.0, 34.0, -1.0)
    18 = @(157.0, -1.0, 32.0, 2.0, 3.0, 36.0, 35.0, 20.0, 29.0, 21.0, 16.0, 11.0, -1.0)
    19 = @(151.0, -1.0, 6.0, 43.0, 7.0, 23.0, 14.0, 25.0, 18.0, -1.0)
    20 = @(159.0, -1.0, 32.0, 2.0, 20.0, 16.0, 50.0, 34.0, 30.0, 11.0, -1.0)
    21 = @(125.0, -1.0, 12.0, 18.0, 41.0, 25.0, -1.0)
}

# The values above start at column D; write them out cell by cell.
$startCol = 4  # column D

foreach ($r in ($rowData.Keys | Sort-Object)) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $col = $startCol + $i
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
